$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G4").Value = 1.45
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 2.05
$ws.Range("AD4").Value = 5.5
$ws.Range("AP4").Value = 26
$ws.Range("AS4").Value = 81
$ws.Range("G5").Value = 2.6
$ws.Range("I5").Value = 2.8
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 7.5
$ws.Range("AD5").Value = 12
$ws.Range("G8").Value = 3.75
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 1.93
$ws.Range("J8").Value = 4.1
$ws.Range("L8").Value = 2.45
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 1.32
$ws.Range("P8").Value = 3.1
$ws.Range("S8").Value = 1.93
$ws.Range("T8").Value = 1.78
$ws.Range("W8").Value = 3.25
$ws.Range("X8").Value = 1.3
$ws.Range("Z8").Value = 2.85
$ws.Range("AA8").Value = 1.8
$ws.Range("AB8").Value = 1.91
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 20
$ws.Range("AE8").Value = 12.5
$ws.Range("AF8").Value = 55
$ws.Range("AG8").Value = 35
$ws.Range("AH8").Value = 40
$ws.Range("AI8").Value = 7
$ws.Range("AJ8").Value = 6.5
$ws.Range("AK8").Value = 14.5
$ws.Range("AL8").Value = 70
$ws.Range("AM8").Value = 600
$ws.Range("AN8").Value = 7
$ws.Range("AO8").Value = 9
$ws.Range("AP8").Value = 8.5
$ws.Range("AR8").Value = 15.5
$ws.Range("AS8").Value = 27
